# Add a new "20201012" sheet to the photo-insert-SQL generator workbook,
# populated the same way as the other date sheets: ID / Shop ID / SQL
# columns, with the SQL column built from a CONCAT() formula.

$wb = $excel.ActiveWorkbook

# --- Deselect / re-anchor the previously active sheet (20201008) -------
# Before the edit it was the active tab, scrolled down with B81 selected.
# After the edit it is just a normal background sheet selected at A1:C2.
$ws8 = $wb.Worksheets.Item("20201008")
$ws8.Activate()
$ws8.Range("A1:C2").Select()

# --- Create the new sheet, positioned after the last existing sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "20201012"

# --- Header row ----------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("A1").Font.Color = 0
$ws.Range("B1").Value = "Shop ID"
$ws.Range("B1").Font.Color = 0
$ws.Range("C1").Value = "SQL"
$ws.Range("C1").Font.Color = 0

# --- Data rows: IDs 13..25, all for the same shop/restaurant uuid ------
$shopId = "da04f9c7-ffb0-11ea-ba65-065a10bcba76"
$ids = 13..25
$row = 2
foreach ($id in $ids) {
    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)

    $aCell.Value = $id
    # Matches the alternating-style pattern used on every other sheet in
    # this workbook: even data rows get the explicit black-font style,
    # odd ones keep the default.
    if ($row % 2 -eq 0) {
        $aCell.Font.Color = 0
    }

    $bCell.Value = $shopId
    $bCell.Font.Color = 0

    $cCell.Formula = "=_xlfn.CONCAT(""INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin('"", B$row, ""'), LPAD("", A$row, "", 7, '0'), 'dish'"", "");"")"

    $row = $row + 1
}

# --- Sheet view: this new sheet becomes the active / selected tab ------
$ws.Activate()
$ws.Range("C2:C14").Select()
